# Scheduled runner update: refresh market-price derived columns (H-N)
# across the Leve profit tables on each class sheet.
$wb = $excel.ActiveWorkbook

$wsALC = $wb.Worksheets.Item("ALC")
$wsALC.Range("H74").Value = 4665
$wsALC.Range("I74").Value = 4490
$wsALC.Range("J74").Value = 4700
$wsALC.Range("K74").Value = 4490
$wsALC.Range("L74").Value = 4700
$wsALC.Range("M74").Value = -3554
$wsALC.Range("N74").Value = -6572
$wsALC.Range("H77").Value = 4665
$wsALC.Range("I77").Value = 4490
$wsALC.Range("J77").Value = 4700
$wsALC.Range("K77").Value = 22450
$wsALC.Range("L77").Value = 23500
$wsALC.Range("M77").Value = -17770
$wsALC.Range("N77").Value = -32860
$wsALC.Range("H132").Value = 34510.035
$wsALC.Range("I132").Value = 44416.914
$wsALC.Range("J132").Value = 1958.8572
$wsALC.Range("K132").Value = 133250.742
$wsALC.Range("L132").Value = 5876.571599999999
$wsALC.Range("M132").Value = -130720.742
$wsALC.Range("N132").Value = -10936.5716
$wsALC.Range("H137").Value = 1814.6786
$wsALC.Range("I137").Value = 1361.6111
$wsALC.Range("J137").Value = 2630.2
$wsALC.Range("K137").Value = 4084.8333
$wsALC.Range("L137").Value = 7890.599999999999
$wsALC.Range("M137").Value = -1534.8333
$wsALC.Range("N137").Value = -12990.6

$wsARM = $wb.Worksheets.Item("ARM")
$wsARM.Range("H32").Value = 4354.9165
$wsARM.Range("I32").Value = 2769.8616
$wsARM.Range("J32").Value = 19073.285
$wsARM.Range("K32").Value = 2769.8616
$wsARM.Range("L32").Value = 19073.285
$wsARM.Range("M32").Value = -2482.8616
$wsARM.Range("N32").Value = -19647.285
$wsARM.Range("H61").Value = 1738.027
$wsARM.Range("I61").Value = 1220.3823
$wsARM.Range("J61").Value = 7604.6665
$wsARM.Range("K61").Value = 1220.3823
$wsARM.Range("L61").Value = 7604.6665
$wsARM.Range("M61").Value = -1008.3823
$wsARM.Range("N61").Value = -8028.6665
$wsARM.Range("H74").Value = 1157.4
$wsARM.Range("I74").Value = 1332
$wsARM.Range("J74").Value = 833.1429000000001
$wsARM.Range("K74").Value = 1332
$wsARM.Range("L74").Value = 833.1429000000001
$wsARM.Range("M74").Value = -458
$wsARM.Range("N74").Value = -2581.1429
$wsARM.Range("H77").Value = 1157.4
$wsARM.Range("I77").Value = 1332
$wsARM.Range("J77").Value = 833.1429000000001
$wsARM.Range("K77").Value = 6660
$wsARM.Range("L77").Value = 4165.7145
$wsARM.Range("M77").Value = -2292
$wsARM.Range("N77").Value = -12901.7145
$wsARM.Range("H132").Value = 2737
$wsARM.Range("I132").Value = 1246
$wsARM.Range("K132").Value = 3738
$wsARM.Range("M132").Value = -1208
$wsARM.Range("H136").Value = 1738.027
$wsARM.Range("I136").Value = 1220.3823
$wsARM.Range("J136").Value = 7604.6665
$wsARM.Range("K136").Value = 3661.1469
$wsARM.Range("L136").Value = 22813.9995
$wsARM.Range("M136").Value = -1111.1469
$wsARM.Range("N136").Value = -27913.9995

$wsBSM = $wb.Worksheets.Item("BSM")
$wsBSM.Range("H115").Value = 39950
$wsBSM.Range("J115").Value = 39950
$wsBSM.Range("L115").Value = 39950
$wsBSM.Range("N115").Value = -43084
$wsBSM.Range("H134").Value = 2148.2222
$wsBSM.Range("I134").Value = 1050.3077
$wsBSM.Range("J134").Value = 5002.8
$wsBSM.Range("K134").Value = 3150.9231
$wsBSM.Range("L134").Value = 15008.4
$wsBSM.Range("M134").Value = -615.9231
$wsBSM.Range("N134").Value = -20078.4

$wsCRP = $wb.Worksheets.Item("CRP")
$wsCRP.Range("H23").Value = 70010
$wsCRP.Range("I23").Value = 0
$wsCRP.Range("K23").Value = 0
$wsCRP.Range("M23").ClearContents()
$wsCRP.Range("H27").Value = 70010
$wsCRP.Range("I27").Value = 0
$wsCRP.Range("K27").Value = 0
$wsCRP.Range("M27").ClearContents()
$wsCRP.Range("H31").Value = 1824.75
$wsCRP.Range("I31").Value = 1676.4857
$wsCRP.Range("J31").Value = 7014
$wsCRP.Range("K31").Value = 1676.4857
$wsCRP.Range("L31").Value = 7014
$wsCRP.Range("M31").Value = -1381.4857
$wsCRP.Range("N31").Value = -7604
$wsCRP.Range("H34").Value = 1824.75
$wsCRP.Range("I34").Value = 1676.4857
$wsCRP.Range("J34").Value = 7014
$wsCRP.Range("K34").Value = 1676.4857
$wsCRP.Range("L34").Value = 7014
$wsCRP.Range("M34").Value = -1474.4857
$wsCRP.Range("N34").Value = -7418
$wsCRP.Range("H58").Value = 1678
$wsCRP.Range("I58").Value = 1158.0476
$wsCRP.Range("J58").Value = 2067.9644
$wsCRP.Range("K58").Value = 1158.0476
$wsCRP.Range("L58").Value = 2067.9644
$wsCRP.Range("M58").Value = -955.0476000000001
$wsCRP.Range("N58").Value = -2473.9644
$wsCRP.Range("H132").Value = 1733.037
$wsCRP.Range("I132").Value = 1161.0555
$wsCRP.Range("J132").Value = 2877
$wsCRP.Range("K132").Value = 3483.1665
$wsCRP.Range("L132").Value = 8631
$wsCRP.Range("M132").Value = -953.1664999999998
$wsCRP.Range("N132").Value = -13691
$wsCRP.Range("H134").Value = 3085.1765
$wsCRP.Range("I134").Value = 3193.6667
$wsCRP.Range("J134").Value = 2824.8
$wsCRP.Range("K134").Value = 9581.000100000001
$wsCRP.Range("L134").Value = 8474.400000000001
$wsCRP.Range("M134").Value = -7046.000100000001
$wsCRP.Range("N134").Value = -13544.4
$wsCRP.Range("H136").Value = 1678
$wsCRP.Range("I136").Value = 1158.0476
$wsCRP.Range("J136").Value = 2067.9644
$wsCRP.Range("K136").Value = 3474.142800000001
$wsCRP.Range("L136").Value = 6203.8932
$wsCRP.Range("M136").Value = -924.1428000000005
$wsCRP.Range("N136").Value = -11303.8932

$wsGSM = $wb.Worksheets.Item("GSM")
$wsGSM.Range("H132").Value = 3171.375
$wsGSM.Range("I132").Value = 2060.389
$wsGSM.Range("J132").Value = 6504.3335
$wsGSM.Range("K132").Value = 6181.167
$wsGSM.Range("L132").Value = 19513.0005
$wsGSM.Range("M132").Value = -3651.167
$wsGSM.Range("N132").Value = -24573.0005

$wsLTW = $wb.Worksheets.Item("LTW")
$wsLTW.Range("H122").Value = 5294862
$wsLTW.Range("I122").Value = 18519952
$wsLTW.Range("K122").Value = 55559856
$wsLTW.Range("M122").Value = -55557406
$wsLTW.Range("H132").Value = 10042.849
$wsLTW.Range("I132").Value = 13462.947
$wsLTW.Range("J132").Value = 5401.2856
$wsLTW.Range("K132").Value = 40388.841
$wsLTW.Range("L132").Value = 16203.8568
$wsLTW.Range("M132").Value = -37858.841
$wsLTW.Range("N132").Value = -21263.8568
$wsLTW.Range("H136").Value = 15877052
$wsLTW.Range("I136").Value = 3780.6667
$wsLTW.Range("J136").Value = 37041416
$wsLTW.Range("K136").Value = 11342.0001
$wsLTW.Range("L136").Value = 111124248
$wsLTW.Range("M136").Value = -8792.000100000001
$wsLTW.Range("N136").Value = -111129348

$wsWVR = $wb.Worksheets.Item("WVR")
$wsWVR.Range("H81").Value = 981.46155
$wsWVR.Range("I81").Value = 719.875
$wsWVR.Range("K81").Value = 1439.75
$wsWVR.Range("M81").Value = -378.75
$wsWVR.Range("H84").Value = 981.46155
$wsWVR.Range("I84").Value = 719.875
$wsWVR.Range("K84").Value = 7198.75
$wsWVR.Range("M84").Value = -1894.75
$wsWVR.Range("H111").Value = 23072
$wsWVR.Range("J111").Value = 23072
$wsWVR.Range("L111").Value = 23072
$wsWVR.Range("N111").Value = -31252
$wsWVR.Range("H132").Value = 3019.5
$wsWVR.Range("I132").Value = 2542.4
$wsWVR.Range("J132").Value = 4041.8572
$wsWVR.Range("K132").Value = 7627.200000000001
$wsWVR.Range("L132").Value = 12125.5716
$wsWVR.Range("M132").Value = -5097.200000000001
$wsWVR.Range("N132").Value = -17185.5716
$wsWVR.Range("H136").Value = 2228.5366
$wsWVR.Range("I136").Value = 790.4286
$wsWVR.Range("J136").Value = 3738.55
$wsWVR.Range("K136").Value = 2371.2858
$wsWVR.Range("L136").Value = 11215.65
$wsWVR.Range("M136").Value = 178.7142000000003
$wsWVR.Range("N136").Value = -16315.65
